$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.267.92'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = '1.868.07'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.89%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.06'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.012'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4665'
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2818'
$ws.Range('E8').Value = '  +2.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06387'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.15'
$ws.Range('E10').Value = '  +2.90%  '
$ws.Range('D11').Value = '1.912.65'
$ws.Range('E11').Value = '  +2.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07604'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.49'
$ws.Range('E13').Value = '  +12.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.976'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6386'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '289.13'
$ws.Range('E16').Value = '  +18.49%  '
$ws.Range('D17').Value = '30.263.40'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.006'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.73'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.134.62'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007351'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.023'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.005'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.34'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.081'
$ws.Range('E26').Value = '  -3.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.31'
$ws.Range('E27').Value = '  +6.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.916'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1083'
$ws.Range('E29').Value = '  +6.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.349'
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.025'
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.797'
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04921'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7203'
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.116'
$ws.Range('E35').Value = '  -3.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.737'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01925'
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.688'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.975'
$ws.Range('E39').Value = '  -1.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8662'
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '105.65'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.013'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.566'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4047'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.32'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.071'
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.830'
$ws.Range('E47').Value = '  +1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1186'
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.86'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05586'
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3716'
$ws.Range('E51').Value = '  +0.33%  '
